$d = $word.ActiveDocument

# --- Change 1: update the "Context Concept Kind / ResourceOccurrence" line ---
$d.Content.Find.Execute(
    "Context Concept Kind / ResourceOccurrence, RHS ResourceOccurrence",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Context Concept Kind / Context ResourceOccurrence Kind, RHS ResourceOccurrence",
    2)

# --- Change 2: add two new list items after
#     "(Amor, Pedro, amante, (Pedro, amaA, Maria));" ---
$idx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*(Amor, Pedro, amante, (Pedro, amaA, Maria));*") {
        $idx = $i
    }
}

# Insert first new paragraph (re-fetch the paragraph from the live
# collection each time so Range.Start/End reflect the mutated document).
$para = $d.Paragraphs.Item($idx)
$para.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Item($idx + 1)
$newPara1.Range.InsertBefore("(Amor, Maria, Pedro, amada);")

# Insert second new paragraph after the first new one.
$newPara1 = $d.Paragraphs.Item($idx + 1)
$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($idx + 2)
$newPara2.Range.InsertBefore("(Empleo, Maria, Pedro, compañera);")
